$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new date column before column O (pushing the "26" date, and its
# formatting, one column to the right into P) and fill the new date (23).
# This also auto-extends the "Noviembre" merge (H1:O1 -> H1:P1).
$ws.Range("O1").EntireColumn.Insert()
$ws.Range("O2").Value = 23

# Minor column width tweaks (M, O, P) to match the new layout.
$ws.Columns.Item(13).ColumnWidth = 4.333333333333334
$ws.Columns.Item(15).ColumnWidth = 3.6666666666666665
$ws.Columns.Item(16).ColumnWidth = 3.4999999999999996

# Update the view: zoom 145% -> 130%, and select E1:P6 instead of E8:H12.
$excel.ActiveWindow.Zoom = 130
$ws.Range("E1:P6").Select() | Out-Null
